$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Npc")

# Column I becomes "skill_value" with the actual per-npc skill identifier
$ws.Cells.Item(1, 9).Value = "skill_value"
$ws.Cells.Item(2, 9).Value = "{(19_DeadlyAttack)}"
$ws.Cells.Item(3, 9).Value = "{(19_CriticalHit)}"
$ws.Cells.Item(4, 9).Value = "{(19_Badbomb)}"

# Column J takes over what used to be column L (drop_table), carrying its fill style
$ws.Cells.Item(1, 10).Value = "drop_table"
$ws.Cells.Item(2, 10).Value = "{(19101, 100)}"
$ws.Cells.Item(3, 10).Value = "{(19101, 100)}"
$ws.Cells.Item(4, 10).Value = "{(19101, 100)}"
$ws.Range("J1:J4").Interior.Color = 65535

# Remove the now-obsolete columns K and L (old enhancer/cooler tube + duplicate drop_table)
$ws.Range("K1:L4").Delete()

$ws.Range("I10").Select()
